{"js": "// \"01-\" is prepended to the title run, the \"_GoBack\" bookmark (previously\n// sitting after \".PEGAR O CREME DENTAL E ABRI\") is moved so it now sits\n// between the new \"01-\" run and the \"ALGORITMO DE COMO ESCOVAR OS DENTES\"\n// run of the title paragraph.\n\nconst body = context.document.body;\n\n// 1) Drop the existing \"_GoBack\" bookmark from wherever it currently is\n//    (after \".PEGAR O CREME DENTAL E ABRI\"). Word only ever keeps a single\n//    \"_GoBack\" bookmark, so remove the old one before re-adding it.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the title text and prepend \"01-\" to it. Word merges text\n//    inserted immediately \"before\" an existing range into the same run,\n//    inheriting that run's formatting (rStyle, color, w14 text effects\u2026).\nconst titleResults = body.search(\"ALGORITMO DE COMO ESCOVAR OS DENTES\", {\n  matchCase: true,\n  matchWholeWord: true\n});\ntitleResults.load(\"text\");\nawait context.sync();\n\nconst titleRange = titleResults.items[0];\ntitleRange.insertText(\"01-\", Word.InsertLocation.before);\nawait context.sync();\n\n// 3) Re-insert the \"_GoBack\" bookmark right after the new \"01-\" text and\n//    before \"ALGORITMO DE COMO ESCOVAR OS DENTES\". Inserting a bookmark in\n//    the middle of a run naturally splits it into two runs (one ending in\n//    \"01-\", one starting with \"ALGORITMO\u2026\"), matching the target markup.\nconst prefixResults = body.search(\"01-\", { matchCase: true });\nprefixResults.load(\"text\");\nawait context.sync();\n\nconst prefixRange = prefixResults.items[0];\nconst afterPrefix = prefixRange.getRange(\"After\");\nafterPrefix.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"01-\" is prepended to the title run, and the \"_GoBack\" bookmark\n# (previously sitting right after \".PEGAR O CREME DENTAL E ABRI\") is moved\n# so it now sits between the new \"01-\" run and the\n# \"ALGORITMO DE COMO ESCOVAR OS DENTES\" run of the title paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark from wherever it currently is\n#    (after \".PEGAR O CREME DENTAL E ABRI\"). Word only ever keeps a single\n#    \"_GoBack\" bookmark, so remove the old one before re-adding it.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the title text and prepend \"01-\" to it. InsertBefore extends\n#    the existing run, inheriting its formatting (rStyle, color, w14 text\n#    effects like shadow/reflection/outline\u2026).\n$titleRange = $d.Content\n$titleRange.Find.Execute(\"ALGORITMO DE COMO ESCOVAR OS DENTES\") | Out-Null\n$titleRange.InsertBefore(\"01-\")\n\n# 3) Re-insert the \"_GoBack\" bookmark right after the new \"01-\" text and\n#    before \"ALGORITMO DE COMO ESCOVAR OS DENTES\". Adding a bookmark in the\n#    middle of a run splits it into two runs (one ending in \"01-\", one\n#    starting with \"ALGORITMO\u2026\"), matching the target markup.\n$prefixRange = $d.Content\n$prefixRange.Find.Execute(\"01-\") | Out-Null\n$splitPoint = $d.Range($prefixRange.End, $prefixRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $splitPoint)\n"}
